$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preserving the original
# "General"-style (no explicit style index), even when the text looks
# like a number (e.g. "546.24") or like a multi-dot price ("58.423.48").
# Forcing NumberFormat="@" before the write stops Excel from silently
# re-typing the cell as a Number; resetting Style back to "Normal"
# afterwards drops the temporary format so the cell XF matches the
# untouched cells around it (no leftover s="..." attribute).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "58.423.48"
Set-TextValue $ws.Cells.Item(2, 5) "  -2.33%  "

Set-TextValue $ws.Cells.Item(3, 4) "2.575.13"
Set-TextValue $ws.Cells.Item(3, 5) "  -3.05%  "

Set-TextValue $ws.Cells.Item(4, 5) "  +0.02%  "

Set-TextValue $ws.Cells.Item(5, 4) "546.24"
Set-TextValue $ws.Cells.Item(5, 5) "  +1.37%  "

Set-TextValue $ws.Cells.Item(6, 4) "143.69"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.66%  "

Set-TextValue $ws.Cells.Item(7, 5) "  +0.11%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.584"
Set-TextValue $ws.Cells.Item(8, 5) "  +1.59%  "

Set-TextValue $ws.Cells.Item(9, 4) "6.74"
Set-TextValue $ws.Cells.Item(9, 5) "  +0.97%  "

Set-TextValue $ws.Cells.Item(10, 4) "0.0998"
Set-TextValue $ws.Cells.Item(10, 5) "  -3.65%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.140"
Set-TextValue $ws.Cells.Item(11, 5) "  +3.75%  "

Set-TextValue $ws.Cells.Item(12, 4) "0.332"
Set-TextValue $ws.Cells.Item(12, 5) "  -2.07%  "

Set-TextValue $ws.Cells.Item(13, 4) "3.031.97"
Set-TextValue $ws.Cells.Item(13, 5) "  -3.28%  "

Set-TextValue $ws.Cells.Item(14, 4) "58.341.07"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.34%  "

Set-TextValue $ws.Cells.Item(15, 4) "20.56"
Set-TextValue $ws.Cells.Item(15, 5) "  -3.24%  "

Set-TextValue $ws.Cells.Item(16, 4) "2.570.50"
Set-TextValue $ws.Cells.Item(16, 5) "  -3.51%  "

Set-TextValue $ws.Cells.Item(17, 5) "  -3.16%  "

Set-TextValue $ws.Cells.Item(18, 4) "4.44"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.31%  "

Set-TextValue $ws.Cells.Item(19, 4) "333.72"
Set-TextValue $ws.Cells.Item(19, 5) "  -3.09%  "

Set-TextValue $ws.Cells.Item(20, 4) "10.00"
Set-TextValue $ws.Cells.Item(20, 5) "  -4.19%  "

Set-TextValue $ws.Cells.Item(21, 4) "6.06"
Set-TextValue $ws.Cells.Item(21, 5) "  -4.52%  "

Set-TextValue $ws.Cells.Item(22, 4) "0.999"
Set-TextValue $ws.Cells.Item(22, 5) "  -0.05%  "

Set-TextValue $ws.Cells.Item(23, 5) "  -0.85%  "

Set-TextValue $ws.Cells.Item(24, 4) "0.422"
Set-TextValue $ws.Cells.Item(24, 5) "  +1.10%  "

Set-TextValue $ws.Cells.Item(25, 4) "0.999"
Set-TextValue $ws.Cells.Item(25, 5) "  +0.00%  "

Set-TextValue $ws.Cells.Item(26, 4) "0.157"
Set-TextValue $ws.Cells.Item(26, 5) "  -5.39%  "

Set-TextValue $ws.Cells.Item(27, 4) "7.07"
Set-TextValue $ws.Cells.Item(27, 5) "  -3.69%  "

Set-TextValue $ws.Cells.Item(28, 4) "0.0₃0736"
Set-TextValue $ws.Cells.Item(28, 5) "  -3.28%  "

Set-TextValue $ws.Cells.Item(29, 4) "0.999"
Set-TextValue $ws.Cells.Item(29, 5) "  +0.08%  "

Set-TextValue $ws.Cells.Item(30, 5) "  -0.80%  "

Set-TextValue $ws.Cells.Item(31, 4) "5.90"
Set-TextValue $ws.Cells.Item(31, 5) "  +0.47%  "

Set-TextValue $ws.Cells.Item(32, 4) "154.31"
Set-TextValue $ws.Cells.Item(32, 5) "  +2.62%  "

Set-TextValue $ws.Cells.Item(33, 4) "18.82"
Set-TextValue $ws.Cells.Item(33, 5) "  -1.14%  "

Set-TextValue $ws.Cells.Item(34, 4) "3.87"
Set-TextValue $ws.Cells.Item(34, 5) "  -4.44%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.859"
Set-TextValue $ws.Cells.Item(35, 5) "  +1.76%  "

Set-TextValue $ws.Cells.Item(36, 4) "1.09"
Set-TextValue $ws.Cells.Item(36, 5) "  -5.35%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.817"
Set-TextValue $ws.Cells.Item(37, 5) "  -3.58%  "

Set-TextValue $ws.Cells.Item(38, 5) "  -3.97%  "

Set-TextValue $ws.Cells.Item(39, 4) "3.56"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.27%  "

Set-TextValue $ws.Cells.Item(40, 4) "279.46"
Set-TextValue $ws.Cells.Item(40, 5) "  -5.18%  "

Set-TextValue $ws.Cells.Item(41, 4) "0.998"
Set-TextValue $ws.Cells.Item(41, 5) "  +0.05%  "

Set-TextValue $ws.Cells.Item(42, 4) "0.591"
Set-TextValue $ws.Cells.Item(42, 5) "  -3.02%  "

Set-TextValue $ws.Cells.Item(43, 4) "10.64"
Set-TextValue $ws.Cells.Item(43, 5) "  -0.76%  "

Set-TextValue $ws.Cells.Item(44, 5) "  -1.22%  "

Set-TextValue $ws.Cells.Item(45, 5) "  -2.88%  "

Set-TextValue $ws.Cells.Item(46, 4) "18.46"
Set-TextValue $ws.Cells.Item(46, 5) "  -5.92%  "

Set-TextValue $ws.Cells.Item(47, 4) "0.0227"
Set-TextValue $ws.Cells.Item(47, 5) "  -0.54%  "

Set-TextValue $ws.Cells.Item(48, 4) "1.896.11"
Set-TextValue $ws.Cells.Item(48, 5) "  -4.47%  "

Set-TextValue $ws.Cells.Item(49, 2) "RenderToken"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Cells.Item(49, 4) "4.40"
Set-TextValue $ws.Cells.Item(49, 5) "  -4.18%  "

Set-TextValue $ws.Cells.Item(50, 2) "InjectiveProtocol"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Cells.Item(50, 4) "17.68"
Set-TextValue $ws.Cells.Item(50, 5) "  -4.45%  "

Set-TextValue $ws.Cells.Item(51, 4) "111.73"
Set-TextValue $ws.Cells.Item(51, 5) "  +1.12%  "
